$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update values in column C for rows 2-4 from "yes" to "no"
$ws.Range("C2").Value = "no"
$ws.Range("C3").Value = "no"
$ws.Range("C4").Value = "no"

# Update the active cell selection on Sheet1 to C5
$ws.Activate()
$ws.Range("C5").Select()
